$wb = $excel.ActiveWorkbook

# ALC row 80
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H80").Value = 424.4737
$ws.Range("I80").Value = 287.5
$ws.Range("K80").Value = 862.5
$ws.Range("M80").Value = 135.5

# ALC row 83
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H83").Value = 424.4737
$ws.Range("I83").Value = 287.5
$ws.Range("K83").Value = 2587.5
$ws.Range("M83").Value = 2404.5

# ALC row 125
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H125").Value = 6724
$ws.Range("I125").Value = 1200
$ws.Range("J125").Value = 7337.778
$ws.Range("K125").Value = 10800
$ws.Range("L125").Value = 66040.00200000001
$ws.Range("M125").Value = -8340
$ws.Range("N125").Value = -70960.00200000001

# ALC row 137
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H137").Value = 305932.8
$ws.Range("I137").Value = 1253.9333
$ws.Range("K137").Value = 3761.7999
$ws.Range("M137").Value = -1211.7999

# ARM row 61
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 35067.332
$ws.Range("I61").Value = 3672.3845
$ws.Range("K61").Value = 3672.3845
$ws.Range("M61").Value = -3460.3845

# ARM row 74
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H74").Value = 102109.336
$ws.Range("I74").Value = 101435.445
$ws.Range("K74").Value = 101435.445
$ws.Range("M74").Value = -100561.445

# ARM row 77
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H77").Value = 102109.336
$ws.Range("I77").Value = 101435.445
$ws.Range("K77").Value = 507177.225
$ws.Range("M77").Value = -502809.225

# ARM row 132
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H132").Value = 15071.956
$ws.Range("I132").Value = 16153.096
$ws.Range("J132").Value = 3720
$ws.Range("K132").Value = 48459.288
$ws.Range("L132").Value = 11160
$ws.Range("M132").Value = -45929.288
$ws.Range("N132").Value = -16220

# ARM row 136
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H136").Value = 35067.332
$ws.Range("I136").Value = 3672.3845
$ws.Range("K136").Value = 11017.1535
$ws.Range("M136").Value = -8467.1535

# BSM row 99
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H99").Value = 2624.5
$ws.Range("I99").Value = 2249.9412
$ws.Range("J99").Value = 3898
$ws.Range("K99").Value = 2249.9412
$ws.Range("L99").Value = 3898
$ws.Range("M99").Value = -751.9412000000002
$ws.Range("N99").Value = -6894

# BSM row 105
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H105").Value = 2488.6667
$ws.Range("I105").Value = 1735
$ws.Range("K105").Value = 1735
$ws.Range("M105").Value = 12

# BSM row 129
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H129").Value = 68950.5
$ws.Range("J129").Value = 68950.5
$ws.Range("L129").Value = 68950.5
$ws.Range("N129").Value = -78950.5

# BSM row 134
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 7711.636
$ws.Range("I134").Value = 7482.8
$ws.Range("K134").Value = 22448.4
$ws.Range("M134").Value = -19913.4

# CRP row 31
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 1821.6052
$ws.Range("J31").Value = 3841
$ws.Range("L31").Value = 3841
$ws.Range("N31").Value = -4431

# CRP row 34
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 1821.6052
$ws.Range("J34").Value = 3841
$ws.Range("L34").Value = 3841
$ws.Range("N34").Value = -4245

# CRP row 58
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H58").Value = 8617.576999999999
$ws.Range("I58").Value = 4567.4736
$ws.Range("J58").Value = 19610.715
$ws.Range("K58").Value = 4567.4736
$ws.Range("L58").Value = 19610.715
$ws.Range("M58").Value = -4364.4736
$ws.Range("N58").Value = -20016.715

# CRP row 99
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H99").Value = 7925.778
$ws.Range("I99").Value = 6224
$ws.Range("K99").Value = 6224
$ws.Range("M99").Value = -4726

# CRP row 107
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H107").Value = 34928.645
$ws.Range("I107").Value = 46918
$ws.Range("K107").Value = 46918
$ws.Range("M107").Value = -44998

# CRP row 126
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H126").Value = 7925.778
$ws.Range("I126").Value = 6224
$ws.Range("K126").Value = 18672
$ws.Range("M126").Value = -16202

# CRP row 132
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H132").Value = 6243.0356
$ws.Range("I132").Value = 3138.7273
$ws.Range("K132").Value = 9416.1819
$ws.Range("M132").Value = -6886.1819

# CRP row 136
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H136").Value = 8617.576999999999
$ws.Range("I136").Value = 4567.4736
$ws.Range("J136").Value = 19610.715
$ws.Range("K136").Value = 13702.4208
$ws.Range("L136").Value = 58832.145
$ws.Range("M136").Value = -11152.4208
$ws.Range("N136").Value = -63932.145

# CRP row 141
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H141").Value = 143130.39
$ws.Range("I141").Value = 15000
$ws.Range("J141").Value = 148954.5
$ws.Range("K141").Value = 15000
$ws.Range("L141").Value = 148954.5
$ws.Range("M141").Value = -9820
$ws.Range("N141").Value = -159314.5

# CUL row 5
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H5").Value = 888.3077
$ws.Range("I5").Value = 855
$ws.Range("J5").Value = 916.8570999999999
$ws.Range("K5").Value = 2565
$ws.Range("L5").Value = 2750.5713
$ws.Range("M5").Value = -2453
$ws.Range("N5").Value = -2974.5713

# CUL row 37
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H37").Value = 114105.82
$ws.Range("J37").Value = 114105.82
$ws.Range("L37").Value = 342317.46
$ws.Range("N37").Value = -342541.46

# CUL row 113
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H113").Value = 2521.7778
$ws.Range("J113").Value = 2746.0667
$ws.Range("L113").Value = 8238.2001
$ws.Range("N113").Value = -12578.2001

# CUL row 122
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H122").Value = 66207.36
$ws.Range("I122").Value = 566.3333
$ws.Range("J122").Value = 115438.125
$ws.Range("K122").Value = 5096.9997
$ws.Range("L122").Value = 1038943.125
$ws.Range("M122").Value = -2646.9997
$ws.Range("N122").Value = -1043843.125

# CUL row 135
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H135").Value = 888.3077
$ws.Range("I135").Value = 855
$ws.Range("J135").Value = 916.8570999999999
$ws.Range("K135").Value = 7695
$ws.Range("L135").Value = 8251.713899999999
$ws.Range("M135").Value = -5160
$ws.Range("N135").Value = -13321.7139

# GSM row 80
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 11408.417
$ws.Range("I80").Value = 3883.5
$ws.Range("J80").Value = 18933.334
$ws.Range("K80").Value = 3883.5
$ws.Range("L80").Value = 18933.334
$ws.Range("M80").Value = -2885.5
$ws.Range("N80").Value = -20929.334

# GSM row 83
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H83").Value = 11408.417
$ws.Range("I83").Value = 3883.5
$ws.Range("J83").Value = 18933.334
$ws.Range("K83").Value = 19417.5
$ws.Range("L83").Value = 94666.67
$ws.Range("M83").Value = -14425.5
$ws.Range("N83").Value = -104650.67

# LTW row 20
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H20").Value = 8617.591
$ws.Range("I20").Value = 5500
$ws.Range("J20").Value = 11215.583
$ws.Range("K20").Value = 5500
$ws.Range("L20").Value = 11215.583
$ws.Range("M20").Value = -5274
$ws.Range("N20").Value = -11667.583

# LTW row 61
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H61").Value = 4619.8
$ws.Range("I61").Value = 1034.3334
$ws.Range("K61").Value = 1034.3334
$ws.Range("M61").Value = -832.3334

# LTW row 68
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H68").Value = 2986.625
$ws.Range("I68").Value = 1299
$ws.Range("J68").Value = 3999.2
$ws.Range("K68").Value = 1299
$ws.Range("L68").Value = 3999.2
$ws.Range("M68").Value = -550
$ws.Range("N68").Value = -5497.2

# LTW row 71
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H71").Value = 2986.625
$ws.Range("I71").Value = 1299
$ws.Range("J71").Value = 3999.2
$ws.Range("K71").Value = 6495
$ws.Range("L71").Value = 19996
$ws.Range("M71").Value = -2751
$ws.Range("N71").Value = -27484

# LTW row 82
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H82").Value = 1611.3478
$ws.Range("I82").Value = 1710.25
$ws.Range("J82").Value = 1558.6
$ws.Range("K82").Value = 1710.25
$ws.Range("L82").Value = 1558.6
$ws.Range("M82").Value = -1349.25
$ws.Range("N82").Value = -2280.6

# LTW row 85
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H85").Value = 1611.3478
$ws.Range("I85").Value = 1710.25
$ws.Range("J85").Value = 1558.6
$ws.Range("K85").Value = 1710.25
$ws.Range("L85").Value = 1558.6
$ws.Range("M85").Value = -462.25
$ws.Range("N85").Value = -4054.6

# LTW row 113
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H113").Value = 4619.8
$ws.Range("I113").Value = 1034.3334
$ws.Range("K113").Value = 1034.3334
$ws.Range("M113").Value = 1135.6666

# LTW row 131
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H131").Value = 78505
$ws.Range("J131").Value = 78505
$ws.Range("L131").Value = 78505
$ws.Range("N131").Value = -88585

# LTW row 132
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H132").Value = 4703.1763
$ws.Range("I132").Value = 3116.1538
$ws.Range("K132").Value = 9348.4614
$ws.Range("M132").Value = -6818.4614

# WVR row 113
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H113").Value = 1479.48
$ws.Range("I113").Value = 1403.0555
$ws.Range("K113").Value = 4209.166499999999
$ws.Range("M113").Value = -2039.166499999999

# WVR row 122
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H122").Value = 5994.6787
$ws.Range("I122").Value = 5687.72
$ws.Range("K122").Value = 17063.16
$ws.Range("M122").Value = -14613.16
